$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the rows that correspond to entries merged away: car_park (16),
# bike_park (17), and vegetation (23). Deleting bottom-up keeps earlier
# row numbers valid.
$ws.Rows.Item(23).Delete()
$ws.Rows.Item(17).Delete()
$ws.Rows.Item(16).Delete()

# Column A was widened to fit the longest remaining label.
$ws.Columns.Item(1).ColumnWidth = 26.77734375

# Reflect the final selection left after the edit.
$ws.Range("A21:XFD21").Select()
